# Fix runopp_map, use DCOPF if ACOPF failed
#
# The source data in several sheets had numeric columns (VMAX/VMIN on
# TGOV1N, pmax/pmin on PV, pmax on Slack) stored as text (shared strings)
# instead of real numbers. This re-enters the correct numeric values so
# downstream OPF tooling (runopp_map) sees numbers, not text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TGOV1N sheet: columns I (VMAX) and J (VMIN), rows 2-11
# ---------------------------------------------------------------------
$tgov1n = $wb.Worksheets.Item("TGOV1N")

$tgov1nVmax = @{2=1.01; 3=1.05; 4=1.05; 5=1.05; 6=1.05; 7=1.05; 8=1.05; 9=1.05; 10=1.05; 11=1.05}
foreach ($row in 2..11) {
    $tgov1n.Cells.Item($row, 9).Value = $tgov1nVmax[$row]   # column I = VMAX
    $tgov1n.Cells.Item($row, 10).Value = 0                  # column J = VMIN
}

$tgov1n.Activate()
$tgov1n.Range("K14").Select()

# ---------------------------------------------------------------------
# Slack sheet: column K (pmax), row 2
# ---------------------------------------------------------------------
$slack = $wb.Worksheets.Item("Slack")
$slack.Cells.Item(2, 11).Value = 15   # column K = pmax

$slack.Activate()
$slack.Range("I22").Select()

# ---------------------------------------------------------------------
# PV sheet: columns K (pmax) and L (pmin), rows 2-10
# ---------------------------------------------------------------------
$pv = $wb.Worksheets.Item("PV")

$pvPmax = @{2=12; 3=7; 4=8; 5=7; 6=7; 7=8; 8=7; 9=7; 10=10}
foreach ($row in 2..10) {
    $pv.Cells.Item($row, 11).Value = $pvPmax[$row]   # column K = pmax
    $pv.Cells.Item($row, 12).Value = 0                # column L = pmin
}

# PV ends up the active / selected sheet (last one touched, like in the
# authored workbook where activeTab points at PV and it is tabSelected).
$pv.Activate()
$pv.Range("K11").Select()
